$d = $word.ActiveDocument

$replacements = @(
    @("2025-04-07 Monday", "2025-04-08 Tuesday"),
    @("479÷3=159, 2", "177÷4=44, 1"),
    @("897÷4=224, 1", "725÷4=181, 1"),
    @("532÷9=59, 1", "791÷2=395, 1"),
    @("375÷2=187, 1", "154÷3=51, 1"),
    @("253÷5=50, 3", "606÷9=67, 3"),
    @("885÷6=147, 3", "330÷9=36, 6"),
    @("889÷2=444, 1", "718÷5=143, 3"),
    @("837÷5=167, 2", "900÷4=225, 0"),
    @("783÷2=391, 1", "900÷4=225, 0"),
    @("174÷2=87, 0", "216÷2=108, 0"),
    @("390÷2=195, 0", "590÷2=295, 0"),
    @("437÷6=72, 5", "524÷8=65, 4"),
    @("306÷6=51, 0", "360÷5=72, 0"),
    @("390÷8=48, 6", "241÷4=60, 1"),
    @("585÷2=292, 1", "696÷7=99, 3"),
    @("132÷4=33, 0", "444÷2=222, 0"),
    @("217÷7=31, 0", "395÷4=98, 3"),
    @("586÷7=83, 5", "688÷5=137, 3"),
    @("384÷3=128, 0", "215÷7=30, 5"),
    @("995÷8=124, 3", "300÷6=50, 0"),
    @("958÷8=119, 6", "229÷5=45, 4"),
    @("189÷9=21, 0", "606÷7=86, 4"),
    @("341÷7=48, 5", "649÷7=92, 5"),
    @("587÷3=195, 2", "924÷5=184, 4"),
    @("790÷4=197, 2", "842÷2=421, 0")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
